# Updates the crypto price/volume snapshot (and fixes the order of a few
# rows whose coin names/links/values had been swapped) to match the
# latest scrape. Price cells that look like plain numbers ("0.698",
# "244.14", ...) are written with a temporary Text format so Excel
# doesn't silently coerce them to doubles (which would also lose
# trailing zeros such as "11.00"); the style is restored to Normal
# immediately after so no stray formatting is left behind. Cells that
# already contain multiple dots (e.g. "44.236.38") or the padded
# percentage strings in column E are never number-like, so they are set
# directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.236.38"
$ws.Range("E2").Value = "  +2.29%  "

$ws.Range("D3").Value = "2.383.44"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.698"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "77.47"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.42%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +29.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +21.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +20.55%  "

$ws.Range("E14").Value = "  +2.55%  "

$ws.Range("E15").Value = "  +7.89%  "

$ws.Range("D16").Value = "2.737.30"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.932"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +7.57%  "

$ws.Range("D18").Value = "2.380.42"
$ws.Range("E18").Value = "  +0.71%  "

$ws.Range("D19").Value = "44.370.57"
$ws.Range("E19").Value = "  +2.51%  "

$ws.Range("E20").Value = "  +2.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "79.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "259.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.80%  "

$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("E25").Value = "  +4.99%  "

$ws.Range("E26").Value = "  +0.76%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.11%  "

$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +19.41%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.84%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.132"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.75%  "

$ws.Range("E33").Value = "  +7.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0765"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.12%  "

$ws.Range("E37").Value = "  +6.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("E40").Value = "  +9.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.32%  "

$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.199"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +18.83%  "

$ws.Range("E45").Value = "  +5.00%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +15.44%  "

$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.04%  "

$ws.Range("E48").Value = "  +6.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.68%  "
